$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New event rows to append below the existing data (rows 3-5).
$data = @(
    @("15651062476621353453", "15", "10", "2022", "Rojan Haun", "4"),
    @("17017948802259489261", "18", "10", "2022", "David Schelle geben", "3"),
    @("17846138601393754605", "25", "10", "2022", "Vincent Box'n", "10")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        # Prefix with an apostrophe so Excel stores every value (including
        # the long numeric IDs) as text instead of silently coercing it to
        # a number and losing precision.
        $cell.Value = "'" + $rowData[$col - 1]
    }
}

# Drop the "quote prefix" text format Excel applied above so the new cells
# keep the workbook's default (unstyled) look, matching the rest of the
# sheet's data rows.
$ws.Range("A6:F8").ClearFormats()
